$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename commodity "elc_win-CHE" -> "elc_won-CHE" everywhere it appears
#    (this is a shared string used ~120 times in the re_profiles sheet,
#    column K "commodity").
# ---------------------------------------------------------------------------
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Cells.Replace("elc_win-CHE", "elc_won-CHE")

# ---------------------------------------------------------------------------
# 2) Re-shuffle the comma separated timeslice-code lists stored in C13/C14
#    of the ev_charging_uc sheet. G7 (=C14) and G8 (=C13) pick the new
#    values up automatically on recalculation.
# ---------------------------------------------------------------------------
$wsEv = $wb.Worksheets.Item("ev_charging_uc")

$wsEv.Range("C13").Value = "S2c0415h07,S2c0415h13,S2c0415h15,S2c0415h17,S2d0427h11,S6aH7,S1aH5,S1b0205h09,S1b0205h14,S1b0205h15,S2aH4,S2aH5,S2c0415h08,S5aH2,S5aH5,S6aH2,S1aH4,S2aH2,S2d0427h15,S4aH6,S1b0205h17,S2d0427h14,S2d0427h16,S5aH4,S6aH3,S6aH4,S1aH6,S2c0415h12,S1b0205h07,S2c0415h10,S2c0415h16,S2c0415h18,S2d0427h12,S3aH3,S1b0205h18,S2aH7,S2d0427h08,S3aH4,S4aH3,S4aH4,S1b0205h16,S2aH6,S2d0427h09,S3aH5,S5aH3,S1b0205h13,S2d0427h13,S4aH5,S1b0205h12,S2c0415h09,S2c0415h11,S3aH7,S1aH2,S1aH7,S1b0205h11,S2aH3,S2d0427h18,S3aH6,S4aH7,S5aH6,S5aH7,S6aH6,S1aH3,S2d0427h07,S2d0427h10,S2d0427h17,S6aH5,S1b0205h08,S1b0205h10,S2c0415h14,S4aH2,S3aH2"

$wsEv.Range("C14").Value = "S1aH1,S1b0205h04,S1b0205h06,S1b0205h21,S2c0415h02,S2c0415h19,S2d0427h03,S3aH1,S2c0415h20,S2c0415h04,S2d0427h21,S4aH1,S6aH1,S1b0205h03,S2c0415h23,S6aH8,S1b0205h23,S2c0415h06,S2d0427h04,S2d0427h24,S1b0205h02,S2d0427h01,S2d0427h19,S1aH8,S1b0205h22,S2d0427h02,S2d0427h06,S2d0427h23,S5aH1,S1b0205h01,S1b0205h20,S2c0415h03,S2c0415h05,S2c0415h24,S2aH8,S2d0427h20,S1b0205h24,S2c0415h22,S2d0427h05,S2d0427h22,S1b0205h05,S2c0415h21,S4aH8,S1b0205h19,S2aH1,S2c0415h01,S5aH8,S3aH8"

# ---------------------------------------------------------------------------
# 3) Reshuffle the season (M) / hydro-share (N) pairing in re_profiles
#    rows 4-9 -- the set of (season, value) pairs is unchanged, only which
#    row each season sits on is different.
# ---------------------------------------------------------------------------
$wsRe.Range("M4").Value = "S4"
$wsRe.Range("N4").Value = 0.19493167383075638

$wsRe.Range("M5").Value = "S3"
$wsRe.Range("N5").Value = 0.26179508564829657

$wsRe.Range("M6").Value = "S1"
$wsRe.Range("N6").Value = 0.17653172515557836

$wsRe.Range("M7").Value = "S5"
$wsRe.Range("N7").Value = 0.080605632899210883

$wsRe.Range("M8").Value = "S2"
$wsRe.Range("N8").Value = 0.40660807082825429

$wsRe.Range("M9").Value = "S6"
$wsRe.Range("N9").Value = 0.079527811637903387
